$d = $word.ActiveDocument

$pairs = @(
    @("153×7=", "126×7="),
    @("437×3=", "889×7="),
    @("766×4=", "360×9="),
    @("255×7=", "602×7="),
    @("403×3=", "969×4="),
    @("425×5=", "567×5="),
    @("903×4=", "685×5="),
    @("413×7=", "177×9="),
    @("955×7=", "406×3="),
    @("135×6=", "702×9="),
    @("879×7=", "531×4="),
    @("406×4=", "991×2="),
    @("101×2=", "241×4="),
    @("131×2=", "735×3="),
    @("296×9=", "585×5="),
    @("972×9=", "357×9="),
    @("642×5=", "451×7="),
    @("921×7=", "832×9="),
    @("289×3=", "525×6="),
    @("746×3=", "583×3="),
    @("682×5=", "338×2="),
    @("498×5=", "151×9="),
    @("534×9=", "176×5="),
    @("179×3=", "912×4="),
    @("706×4=", "483×2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
